$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the courses already known up to "Introduction to Predictive Analytics in Python"
$part1 = @(
    "Supervised Learning with Scikit Learn",
    "Unsupervised Learning in Python",
    "Machine Learning with Tree-Based Models in Python",
    "Linear Classifiers in Python",
    "Cluster Analysis in Python",
    "Extreme Gradient Boosting with XGBoost",
    "Preprocessing for Machine Learning in Python",
    "Feature Engineering for ML in Python",
    "Dimensionality Reduction in Python",
    "Model Validation in Python",
    "Hypertuning in Python",
    "Introduction to Predictive Analytics in Python"
)

$row = 2
foreach ($course in $part1) {
    $ws.Cells.Item($row, 1).Value = $course
    $row++
}

# Continue with the next courses
$ws.Cells.Item(14, 1).Value = "Ensemble Methods in Python"
$ws.Cells.Item(15, 1).Value = "Practicing ML Interview Questions in Python"
$ws.Cells.Item(16, 1).Value = "Designing ML Workflows in Python"

# Go back and insert a course that was missed, shifting the rows below down
$ws.Rows("14").Insert()
$ws.Cells.Item(14, 1).Value = "Intermediate Predictive Analytics in Python"

# Remaining courses
$part2 = @(
    "Unit Testing in Python",
    "Software Engineering for Data Scientist in Python",
    "Practicing Coding Interview Questions in Python",
    "Assessment: ML Fundamentals in Python",
    "Assessment: Statistics Fundamentals with Python",
    "Assessment: Data Analysis in SQL",
    "Introduction to SQL",
    "Intermediate SQL",
    "Joining Data in SQL",
    "EDA in SQL",
    "Functions for Manipulating Data in PSQL"
)

$row = 18
foreach ($course in $part2) {
    $ws.Cells.Item($row, 1).Value = $course
    $row++
}

# Add a "Completed Date" tracking column
$ws.Range("B1").Value = "Completed Date"
$ws.Range("B1").Font.Bold = $true
$ws.Columns("B").EntireColumn.AutoFit()

# Update the view (zoom level and selected cell)
$ws.Range("D3").Select()
$excel.ActiveWindow.Zoom = 150
